$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The "BRCNN 128 Seed 10" run was actually just "BRCNN 128" (seed noted elsewhere) -
# fix the label in column A for the BRCNN rows (8-13).
foreach ($r in 8..13) {
    $ws.Cells.Item($r, 1).Value = "BRCNN 128"
}

# Corrected f1-score for the "BRCNN 128 / Limited to start / combined" row.
$ws.Range("F13").Value = 0.52

# Hide the Goldstandard rows (2-7) - they're not the focus of this comparison anymore.
$ws.Range("A2:A7").EntireRow.Hidden = $true
